$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ------------------------------------------------------------------
# Insert three new rows: one at position 5 (new "schools" Q&A,
# pushing the old rows 5 & 6 down to 6 & 7), and two more after the
# (now shifted) forest-productivity row, at positions 8 and 9, for
# the new "coniferous/deciduous" and "cultural heritage" Q&As.
# ------------------------------------------------------------------
$ws.Rows("5:5").Insert()
$ws.Rows("8:8").Insert()
$ws.Rows("9:9").Insert()

# ------------------------------------------------------------------
# Row 5 (new): "How many schools buildings are there?"
# ------------------------------------------------------------------
$ws.Range("A5").Value = "How many schools buildings are there?"
$ws.Range("B5").Value = ""
$ws.Range("C5").Value = "Steps:`n1) Gather and load building point  data and its documentation`n2) Find all ""bygningstype"" values that correspond to some type of school building`n3) Count the rows after filtering"

$ws.Range("A3").Copy()
$ws.Range("A5").PasteSpecial($xlPasteFormats)
$ws.Range("B4").Copy()
$ws.Range("B5").PasteSpecial($xlPasteFormats)
$ws.Range("C3").Copy()
$ws.Range("C5").PasteSpecial($xlPasteFormats)

$ws.Rows("5:5").RowHeight = 86.4

# ------------------------------------------------------------------
# Row 6 (was row 5): unchanged content, just shifted down - nothing
# else to do, values/format moved automatically with the insert.
# ------------------------------------------------------------------

# ------------------------------------------------------------------
# Row 7 (was row 6): the "Reasoning" cell gets new wording (two
# trailing periods removed from the AR50 steps text).
# ------------------------------------------------------------------
$ws.Range("C7").Value = "Steps:`n1) Gather and load building point AR50 data , either the entire series or only ""Jordbruk"", along with the AR50 documentation`n2) Find the ""skogbonitet"" field and check the documentation to see that areas with the highest forest productivity have value 18`n3) Select building points that have ""skogbonitet"" = 18`n4) Calculate the area of each selected polygon`n5) Sum all areas"

# ------------------------------------------------------------------
# Row 8 (new): "Is the area predominantly coniferous or deciduous
# forest?"
# ------------------------------------------------------------------
$ws.Range("A8").Value = "Is the area predominantly coniferous or deciduous forest? "
$ws.Range("B8").Value = "There is about 250km^2 coniferous forest and 17km^2 deciduous forest. Thus, the area is predominantly coniferous. "
$ws.Range("C8").Value = "Steps:`n1) Gather and load building point AR50 data , either the entire series or only ""Jordbruk"", along with the AR50 documentation`n2) Find the ""skogbonitet"" field and check the documentation to see that areas with the highest forest productivity have value 18`n3) Select building points that have ""skogbonitet"" = 18`n4) Calculate the area of each selected polygon.`n5) Sum all areas"

$ws.Range("B6").Copy()
$ws.Range("B8").PasteSpecial($xlPasteFormats)
$ws.Range("C6").Copy()
$ws.Range("C8").PasteSpecial($xlPasteFormats)

# A8 uses a new, lighter style: default (no border / no fill / no
# wrap) font, just switched to Calibri 11 explicitly.
$ws.Range("A8").Font.Name = "Calibri"
$ws.Range("A8").Font.Size = 11
$ws.Range("A8").Borders.LineStyle = -4142
$ws.Range("A8").WrapText = $false

$ws.Rows("8:8").RowHeight = 144

# ------------------------------------------------------------------
# Row 9 (new): cultural heritage sites question.
# ------------------------------------------------------------------
$ws.Range("A9").Value = "Are there any cultural heritage sites within this bounding box?`n[`n    [63.4159261840723, 10.449447170396198],`n    [63.42364892527119, 10.449837255689479],`n    [63.42345265401785, 10.46906466674888],`n    [63.4157299786334, 10.46866941273723],`n    [63.4159261840723, 10.449447170396198] `n]"
$ws.Range("B9").Value = "There are three cultural heritage sites within the bounding box. "
$ws.Range("C9").Value = "Steps:`n1) Gather and load building point cultural heritage point data`n2) Create a bounding box polygon`n3) Do ""is-point-within"" analysis, filtering the cultural heritage data`n4) Count the number of rows after filtering"

$ws.Range("A6").Copy()
$ws.Range("A9").PasteSpecial($xlPasteFormats)
$ws.Range("B6").Copy()
$ws.Range("B9").PasteSpecial($xlPasteFormats)
$ws.Range("C6").Copy()
$ws.Range("C9").PasteSpecial($xlPasteFormats)

$ws.Rows("9:9").RowHeight = 129.6

$ws.Application.CutCopyMode = $false

# ------------------------------------------------------------------
# Sheet view: scroll position & selection.
# ------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("B5").Select()
